# Update the cryptos list (Price + Volume(1h) columns) with freshly
# scraped figures, as produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.762.51"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.387.93"
$ws.Range("E3").Value = "  -1.90%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'568.68"
$ws.Range("E5").Value = "  -2.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'141.28"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.388.43"
$ws.Range("E8").Value = "  -1.89%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.33%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'7.48"
$ws.Range("E10").Value = "  -1.88%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.56%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  +2.12%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.964.19"
$ws.Range("E13").Value = "  -1.94%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'28.37"

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.94%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.75%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.386.58"
$ws.Range("E17").Value = "  -2.09%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "60.858.96"
$ws.Range("E18").Value = "  -1.53%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.26%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'14.01"
$ws.Range("E20").Value = "  -2.12%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'9.00"
$ws.Range("E21").Value = "  -5.91%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'383.87"
$ws.Range("E22").Value = "  -1.62%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = "  -1.04%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'73.76"
$ws.Range("E24").Value = "  +0.23%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.22%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -5.62%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "3.525.40"
$ws.Range("E27").Value = "  -1.89%  "

# Row 28 - Kaspa
$ws.Range("D28").Value = "'0.178"
$ws.Range("E28").Value = "  -2.39%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 - RenderToken
$ws.Range("E30").Value = "  -2.66%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "  -2.11%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -2.25%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -2.56%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.02%  "

# Row 35 - EthereumClassic
$ws.Range("D35").Value = "'23.72"
$ws.Range("E35").Value = "  -1.51%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'6.98"
$ws.Range("E36").Value = "  -0.55%  "

# Row 37 - Monero
$ws.Range("D37").Value = "'166.57"
$ws.Range("E37").Value = "  -0.02%  "

# Row 39 - RenzoRestakedETH
$ws.Range("D39").Value = "3.418.96"
$ws.Range("E39").Value = "  -1.77%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "'28.12"
$ws.Range("E41").Value = "  +1.37%  "

# Row 42 - Hedera
$ws.Range("D42").Value = "'0.0775"
$ws.Range("E42").Value = "  -0.97%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -3.01%  "

# Row 44 - FirstDigitalUSD
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.09%  "

# Row 45 - OKB
$ws.Range("D45").Value = "'41.82"
$ws.Range("E45").Value = "  -1.54%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "  -1.87%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -2.97%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.525.93"
$ws.Range("E48").Value = "  -1.90%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  -3.16%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'23.67"
$ws.Range("E50").Value = "  +2.76%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  -1.30%  "
